# "Generate Report for Handoff"
#
# The localization-status report regenerated its handoff Xliff for the
# 0b8bc02d-d3d0-4ca9-807c-8e026ea205e9 source file: the "Latest Handoff
# Datetime" for the zh-cn target moved forward from 02:41:16 to 02:41:32.
#
# Update the zh-cn worksheet's row for that file (row 5), column H
# ("Latest Handoff Datetime"), to reflect the new handoff timestamp.
# The value is stored as plain text (not a real Excel date serial), so
# write it as a literal string to keep the same cell type/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("H5").Value = "2016-08-24 02:41:32"
